$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E11:F70 (skip row 21 which stays anchored at 1050) so that the
# "money" stat values used to populate a "diamond" cost column follow the
# uniform progression: value(row) = 1050 + (row - 21) * 50
for ($row = 11; $row -le 70; $row++) {
    if ($row -eq 21) { continue }
    $value = 1050 + (($row - 21) * 50)
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 6).Value = $value
}

# Move the active selection to G30, matching the saved view state.
$ws.Range("G30").Select()
